$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the typo in the URL for row 5 (double dash -> single dash)
$ws.Range("C5").Value = "https://app.datacamp.com/learn/career-tracks/associate-data-scientist-in-python"

# Add new rows 6-8 with new track entries
$ws.Range("A6").Value = 4
$ws.Range("B6").Value = "Excel Fundamentals"
$ws.Range("C6").Value = "https://app.datacamp.com/learn/career-tracks/excel-fundamentals"
$ws.Range("D6").Value = 0

$ws.Range("A7").Value = 5
$ws.Range("B7").Value = "Git Fundamentals"
$ws.Range("C7").Value = "https://app.datacamp.com/learn/career-tracks/git-fundamentals"
$ws.Range("D7").Value = 0

$ws.Range("A8").Value = 6
$ws.Range("B8").Value = "SQL Fundamentals"
$ws.Range("C8").Value = "https://app.datacamp.com/learn/career-tracks/sql-fundamentals"
$ws.Range("D8").Value = 0
